$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "96.688.08"
$ws.Range("E2").Value = "  +0.39%  "

$ws.Range("D3").Value = "3.647.52"
$ws.Range("E3").Value = "  +1.95%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.72%  "

$ws.Range("E6").Value = "  +19.66%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "655.11"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.17%  "

$ws.Range("E8").Value = "  +4.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.06"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.78%  "

$ws.Range("E10").Value = "  -0.04%  "

$ws.Range("D11").Value = "3.644.92"
$ws.Range("E11").Value = "  +1.93%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "44.28"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.22%  "

$ws.Range("E13").Value = "  +1.13%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.51"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.21%  "

$ws.Range("D15").Value = "4.324.62"
$ws.Range("E15").Value = "  +1.85%  "

$ws.Range("D16").Value = "96.432.71"
$ws.Range("E16").Value = "  +0.17%  "

$ws.Range("E17").Value = "  +0.48%  "

$ws.Range("D18").Value = "3.642.96"
$ws.Range("E18").Value = "  +1.70%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.96"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.10%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.77"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.15%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.57%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.536"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +8.94%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.44"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.05%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "512.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.21%  "

$ws.Range("E25").Value = "  +2.78%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.90"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.98%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "101.15"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.99%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "13.08"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.10%  "

$ws.Range("E29").Value = "  +12.37%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.03"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.88%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.88"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.86%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.05%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.185"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.87%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "33.14"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.06%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.20%  "

$ws.Range("E36").Value = "  +8.79%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.584"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.37%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.82"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.11%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "615.15"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.01%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "42.65"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +25.02%  "

$ws.Range("E41").Value = "  +4.92%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.952"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.18%  "

$ws.Range("E43").Value = "  +6.39%  "

$ws.Range("E44").Value = "  -0.02%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.16"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.70%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0443"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.15%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.31"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.70%  "

$ws.Range("E48").Value = "  +0.34%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.414"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +19.42%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.61"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.26%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "54.25"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.57%  "
